$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.064.50'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.248.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.49%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.244.94'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.537'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.159'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.65'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.489'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.78%  '
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.785.05'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.179.28'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.250.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '496.09'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.741'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.97'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.66'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.66%  '
$ws.Range("E30").Value = '  +43.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.97'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.82'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.65'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  -4.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.35'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.41%  '
$ws.Range("E37").Value = '  +14.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.42'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '489.49'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.78%  '
$ws.Range("E40").Value = '  +4.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0418'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.68'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.981.47'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.288'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.24'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.118'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.95%  '
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.74'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.37%  '
